# Applies the "Updated symbol list" GitHub Actions commit to the crypto price
# table on Sheet1. The sheet stores every data cell as literal text (inline
# strings) -- prices, percentage-change figures, and the "Hora" (hour) code
# are all text that merely *look* numeric, and row 23/24 swap the
# ZBToken / CoinExToken entries (name, link, price, volume all move together).
#
# Because Excel auto-converts numeric-looking text (e.g. "4", "286.81",
# "4.25%") into real numbers (or percentages) when you assign .Value, each
# target cell's NumberFormat is forced to "@" (Text) *before* the assignment
# so the literal string is preserved exactly as in the source data, matching
# the original inline-string cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '286.81' },
    @{ Cell = 'E2'; Value = '4.25%' },
    @{ Cell = 'G2'; Value = '4' },
    @{ Cell = 'D3'; Value = '28.27' },
    @{ Cell = 'E3'; Value = '4.09%' },
    @{ Cell = 'G3'; Value = '4' },
    @{ Cell = 'D4'; Value = '4.920' },
    @{ Cell = 'E4'; Value = '1.08%' },
    @{ Cell = 'G4'; Value = '4' },
    @{ Cell = 'E5'; Value = '2.46%' },
    @{ Cell = 'G5'; Value = '4' },
    @{ Cell = 'D6'; Value = '7.243' },
    @{ Cell = 'E6'; Value = '4.50%' },
    @{ Cell = 'G6'; Value = '4' },
    @{ Cell = 'D7'; Value = '1.351' },
    @{ Cell = 'E7'; Value = '11.74%' },
    @{ Cell = 'G7'; Value = '4' },
    @{ Cell = 'D8'; Value = '0.9165' },
    @{ Cell = 'E8'; Value = '4.52%' },
    @{ Cell = 'G8'; Value = '4' },
    @{ Cell = 'D9'; Value = '0.1569' },
    @{ Cell = 'E9'; Value = '3.66%' },
    @{ Cell = 'G9'; Value = '4' },
    @{ Cell = 'D10'; Value = '0.06538' },
    @{ Cell = 'E10'; Value = '28.19%' },
    @{ Cell = 'G10'; Value = '4' },
    @{ Cell = 'E11'; Value = '1.88%' },
    @{ Cell = 'G11'; Value = '4' },
    @{ Cell = 'D12'; Value = '0.02981' },
    @{ Cell = 'E12'; Value = '0.58%' },
    @{ Cell = 'G12'; Value = '4' },
    @{ Cell = 'D13'; Value = '0.08983' },
    @{ Cell = 'E13'; Value = '-0.02%' },
    @{ Cell = 'G13'; Value = '4' },
    @{ Cell = 'D14'; Value = '0.001587' },
    @{ Cell = 'E14'; Value = '1.31%' },
    @{ Cell = 'G14'; Value = '4' },
    @{ Cell = 'D15'; Value = '0.0006539' },
    @{ Cell = 'E15'; Value = '2.26%' },
    @{ Cell = 'G15'; Value = '4' },
    @{ Cell = 'D16'; Value = '0.006109' },
    @{ Cell = 'E16'; Value = '-1.19%' },
    @{ Cell = 'G16'; Value = '4' },
    @{ Cell = 'D17'; Value = '3.488' },
    @{ Cell = 'G17'; Value = '4' },
    @{ Cell = 'D18'; Value = '3.391' },
    @{ Cell = 'G18'; Value = '4' },
    @{ Cell = 'D19'; Value = '2.239' },
    @{ Cell = 'E19'; Value = '-1.97%' },
    @{ Cell = 'G19'; Value = '4' },
    @{ Cell = 'G20'; Value = '4' },
    @{ Cell = 'D21'; Value = '0.1349' },
    @{ Cell = 'E21'; Value = '-0.04%' },
    @{ Cell = 'G21'; Value = '4' },
    @{ Cell = 'D22'; Value = '3.981' },
    @{ Cell = 'E22'; Value = '1.98%' },
    @{ Cell = 'G22'; Value = '4' },
    @{ Cell = 'B23'; Value = 'CoinExToken' },
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' },
    @{ Cell = 'D23'; Value = '0.04470' },
    @{ Cell = 'E23'; Value = '0.95%' },
    @{ Cell = 'G23'; Value = '4' },
    @{ Cell = 'B24'; Value = 'ZBToken' },
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb' },
    @{ Cell = 'D24'; Value = '0.1519' },
    @{ Cell = 'E24'; Value = '10.08%' },
    @{ Cell = 'G24'; Value = '4' },
    @{ Cell = 'D25'; Value = '0.001188' },
    @{ Cell = 'E25'; Value = '1.12%' },
    @{ Cell = 'G25'; Value = '4' },
    @{ Cell = 'D26'; Value = '0.004340' },
    @{ Cell = 'E26'; Value = '12.54%' },
    @{ Cell = 'G26'; Value = '4' },
    @{ Cell = 'G27'; Value = '4' },
    @{ Cell = 'E28'; Value = '-1.75%' },
    @{ Cell = 'G28'; Value = '4' },
    @{ Cell = 'D29'; Value = '0.0001635' },
    @{ Cell = 'E29'; Value = '-15.74%' },
    @{ Cell = 'G29'; Value = '4' },
    @{ Cell = 'G30'; Value = '4' },
    @{ Cell = 'G31'; Value = '4' },
    @{ Cell = 'G32'; Value = '4' },
    @{ Cell = 'G33'; Value = '4' },
    @{ Cell = 'G34'; Value = '4' },
    @{ Cell = 'G35'; Value = '4' },
    @{ Cell = 'G36'; Value = '4' },
    @{ Cell = 'G37'; Value = '4' },
    @{ Cell = 'G38'; Value = '4' },
    @{ Cell = 'G39'; Value = '4' },
    @{ Cell = 'D40'; Value = '0.04166' },
    @{ Cell = 'E40'; Value = '0.91%' },
    @{ Cell = 'G40'; Value = '4' },
    @{ Cell = 'D41'; Value = '0.007000' },
    @{ Cell = 'E41'; Value = '2.82%' },
    @{ Cell = 'G41'; Value = '4' },
    @{ Cell = 'D42'; Value = '0.1417' },
    @{ Cell = 'E42'; Value = '20.70%' },
    @{ Cell = 'G42'; Value = '4' },
    @{ Cell = 'D43'; Value = '0.002059' },
    @{ Cell = 'E43'; Value = '-1.98%' },
    @{ Cell = 'G43'; Value = '4' },
    @{ Cell = 'D44'; Value = '0.01245' },
    @{ Cell = 'E44'; Value = '8.48%' },
    @{ Cell = 'G44'; Value = '4' },
    @{ Cell = 'D45'; Value = '0.00005549' },
    @{ Cell = 'E45'; Value = '6.99%' },
    @{ Cell = 'G45'; Value = '4' },
    @{ Cell = 'D46'; Value = '1.562' },
    @{ Cell = 'E46'; Value = '-7.84%' },
    @{ Cell = 'G46'; Value = '4' },
    @{ Cell = 'D47'; Value = '0.01850' },
    @{ Cell = 'E47'; Value = '-7.54%' },
    @{ Cell = 'G47'; Value = '4' },
    @{ Cell = 'G48'; Value = '4' },
    @{ Cell = 'G49'; Value = '4' },
    @{ Cell = 'G50'; Value = '4' },
    @{ Cell = 'G51'; Value = '4' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}

Write-Output ("Updated {0} cells on '{1}'" -f $updates.Count, $ws.Name)
